# Application of safety margins for collision avoidance
# Updates target x/y waypoint coordinates (and a couple of derived z values)
# on the active worksheet to reflect newly computed collision-avoidance
# safety margins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 291.0
$ws.Range("B3").Value = 250.0
$ws.Range("A4").Value = 292.0
$ws.Range("B4").Value = 252.0
$ws.Range("A5").Value = 297.0
$ws.Range("B5").Value = 250.0
$ws.Range("C5").Value = 23.8350718518842
$ws.Range("A6").Value = 298.0
$ws.Range("B6").Value = 250.0
$ws.Range("C6").Value = 27.410332629666833
$ws.Range("A7").Value = 299.0
$ws.Range("B7").Value = 255.0
$ws.Range("A8").Value = 301.0
$ws.Range("B8").Value = 245.0
$ws.Range("A9").Value = 304.0
$ws.Range("B9").Value = 244.0
$ws.Range("A10").Value = 304.0
$ws.Range("B10").Value = 256.0
$ws.Range("A11").Value = 304.0
$ws.Range("B11").Value = 259.0
$ws.Range("A12").Value = 304.0
$ws.Range("B12").Value = 261.0
$ws.Range("A13").Value = 306.0
$ws.Range("B13").Value = 264.0
$ws.Range("A14").Value = 309.0
$ws.Range("B14").Value = 266.0
$ws.Range("A15").Value = 309.0
$ws.Range("B15").Value = 269.0
$ws.Range("A16").Value = 310.0
$ws.Range("B16").Value = 272.0
$ws.Range("A17").Value = 314.0
$ws.Range("B17").Value = 274.0
$ws.Range("A18").Value = 319.0
$ws.Range("B18").Value = 273.0
$ws.Range("A19").Value = 321.0
$ws.Range("B19").Value = 275.0
$ws.Range("A20").Value = 321.0
$ws.Range("B20").Value = 281.0
$ws.Range("A21").Value = 326.0
$ws.Range("B21").Value = 281.0
$ws.Range("A22").Value = 325.0
$ws.Range("B22").Value = 286.0
$ws.Range("A23").Value = 326.0
$ws.Range("B23").Value = 290.0
$ws.Range("A24").Value = 328.0
$ws.Range("B24").Value = 293.0
$ws.Range("A25").Value = 329.0
$ws.Range("B25").Value = 296.0
$ws.Range("A26").Value = 324.0
$ws.Range("B26").Value = 299.0
$ws.Range("A27").Value = 334.0
$ws.Range("B27").Value = 299.0
$ws.Range("A28").Value = 334.0
$ws.Range("B28").Value = 305.0
$ws.Range("A29").Value = 334.0
$ws.Range("B29").Value = 309.0
$ws.Range("A30").Value = 335.0
$ws.Range("B30").Value = 312.0
$ws.Range("A31").Value = 336.0
$ws.Range("A32").Value = 339.0
$ws.Range("B32").Value = 319.0
$ws.Range("A33").Value = 344.0
$ws.Range("B33").Value = 319.0
$ws.Range("A34").Value = 348.0
$ws.Range("B34").Value = 318.0
$ws.Range("A35").Value = 350.0
$ws.Range("B35").Value = 318.0
$ws.Range("A36").Value = 350.0
$ws.Range("B36").Value = 327.0
$ws.Range("A37").Value = 354.0
$ws.Range("B37").Value = 329.0
$ws.Range("A38").Value = 354.0
$ws.Range("B38").Value = 333.0
$ws.Range("A39").Value = 356.0
$ws.Range("B39").Value = 336.0
$ws.Range("A40").Value = 362.0
$ws.Range("B40").Value = 335.0
$ws.Range("A41").Value = 364.0
$ws.Range("A42").Value = 364.0
$ws.Range("B42").Value = 344.0
$ws.Range("A43").Value = 364.0
$ws.Range("B43").Value = 348.0
$ws.Range("A44").Value = 369.0
$ws.Range("B44").Value = 349.0
$ws.Range("A45").Value = 370.0
$ws.Range("B45").Value = 353.0
$ws.Range("A46").Value = 375.0
$ws.Range("A47").Value = 374.0
$ws.Range("B47").Value = 359.0
$ws.Range("A48").Value = 375.0
$ws.Range("B48").Value = 363.0
$ws.Range("A49").Value = 381.0
$ws.Range("B49").Value = 364.0
$ws.Range("A50").Value = 381.0
$ws.Range("B50").Value = 369.0
$ws.Range("A51").Value = 386.0
$ws.Range("B51").Value = 370.0
$ws.Range("A52").Value = 386.0
$ws.Range("B52").Value = 375.0
$ws.Range("A53").Value = 388.0
$ws.Range("B53").Value = 378.0
$ws.Range("A54").Value = 389.0
$ws.Range("B54").Value = 381.0
$ws.Range("A55").Value = 389.0
$ws.Range("B55").Value = 386.0
$ws.Range("A56").Value = 394.0
$ws.Range("B56").Value = 388.0
$ws.Range("A57").Value = 399.0
$ws.Range("B57").Value = 388.0
$ws.Range("A58").Value = 402.0
$ws.Range("B58").Value = 390.0
$ws.Range("A59").Value = 404.0
$ws.Range("B59").Value = 394.0
$ws.Range("A60").Value = 404.0
$ws.Range("B60").Value = 399.0
$ws.Range("A61").Value = 409.0
$ws.Range("B61").Value = 399.0
$ws.Range("A62").Value = 413.0
$ws.Range("B62").Value = 400.0
$ws.Range("A63").Value = 417.0
$ws.Range("B63").Value = 400.0
$ws.Range("A64").Value = 420.0
$ws.Range("B64").Value = 400.0
$ws.Range("A65").Value = 424.0
$ws.Range("B65").Value = 401.0
$ws.Range("A66").Value = 425.0
$ws.Range("B66").Value = 395.0
$ws.Range("A67").Value = 429.0
$ws.Range("B67").Value = 406.0
$ws.Range("A68").Value = 432.0
$ws.Range("B68").Value = 408.0
$ws.Range("A69").Value = 435.0
$ws.Range("B69").Value = 410.0
$ws.Range("A70").Value = 430.0
$ws.Range("B70").Value = 415.0
$ws.Range("A71").Value = 435.0
$ws.Range("B71").Value = 417.0
$ws.Range("A72").Value = 435.0
$ws.Range("B72").Value = 421.0
$ws.Range("A73").Value = 436.0
$ws.Range("B73").Value = 424.0
$ws.Range("A74").Value = 440.0
$ws.Range("B74").Value = 427.0
$ws.Range("A75").Value = 436.0
$ws.Range("B75").Value = 430.0
$ws.Range("A76").Value = 443.0
$ws.Range("B76").Value = 433.0
$ws.Range("A77").Value = 445.0
$ws.Range("B77").Value = 436.0
$ws.Range("A78").Value = 445.0
$ws.Range("B78").Value = 440.0
$ws.Range("A79").Value = 450.0
$ws.Range("B79").Value = 441.0
$ws.Range("A80").Value = 454.0
$ws.Range("B80").Value = 441.0
$ws.Range("A81").Value = 455.0
$ws.Range("B81").Value = 446.0
